$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E4").Value = 3711
$ws.Range("F4").Value = 3667
$ws.Range("G4").Value = 3623
$ws.Range("H4").Value = 3630
$ws.Range("I4").Value = 3584
$ws.Range("J4").Value = 3590
$ws.Range("K4").Value = 3543

$ws.Range("E4:K4").Select()
